$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

# Column D is a date, formatted like the other rows (style with numFmtId 165)
$baseDate = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($row, 4).Value = $baseDate.AddDays(44890)
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = 300000000
$ws.Cells.Item($row, 7).Value = "Espárragos"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 160
$ws.Cells.Item($row, 11).Value = 900
$ws.Cells.Item($row, 12).Value = 1000
$ws.Cells.Item($row, 13).Value = 950
$ws.Cells.Item($row, 14).Value = "`$/kilo"
$ws.Cells.Item($row, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($row, 16).Value = 950
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 18).Value = "Hortaliza"
